$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SchemaOrganization")
$ws.Range("B2").Value = 'http://example.com/organization1:Image1'
$ws.Range("B3").Value = 'http://example.com/organization3:Image1'
$ws.Range("B4").Value = 'http://example.com/organization5:Image1'

$ws = $wb.Worksheets.Item("RightsStatementsDotOrgRightsStatement")
$ws.Range("E2").Value = 'You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material.'
$ws.Range("E3").Value = 'Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use.'

$ws = $wb.Worksheets.Item("FoafPerson")
$ws.Range("E2").Value = 'http://example.com/person0:Image0'
$ws.Range("E3").Value = 'http://example.com/person2:Image1'

$ws = $wb.Worksheets.Item("RdfProperty")
$ws.Range("C2").Value = 'dcterms:description:Image0'
$ws.Range("C3").Value = 'dcterms:extent:Image0'
$ws.Range("C4").Value = 'dcterms:language:Image0'
$ws.Range("C6").Value = 'dcterms:publisher:Image1'
$ws.Range("C8").Value = 'dcterms:spatial:Image1'
$ws.Range("C9").Value = 'dcterms:subject:Image1'
$ws.Range("C10").Value = 'dcterms:title:Image1'

$ws = $wb.Worksheets.Item("SchemaProperty")
$ws.Range("C2").Value = 'schema:description:Image1'
$ws.Range("C4").Value = 'schema:spatial:Image0'

$ws = $wb.Worksheets.Item("FoafOrganization")
$ws.Range("C3").Value = 'http://example.com/organization2:Image1'
$ws.Range("C4").Value = 'http://example.com/organization4:Image0'

$ws = $wb.Worksheets.Item("SkosConcept")
$ws.Range("B3").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image0'
$ws.Range("B5").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:3:Image1'
$ws.Range("B6").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:4:Image0'
$ws.Range("B9").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:7:Image1'
$ws.Range("B10").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image0'
$ws.Range("B13").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:11:Image1'
$ws.Range("B16").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image1'
$ws.Range("B17").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image0'
$ws.Range("B19").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:17:Image0'
$ws.Range("B21").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:19:Image0'
$ws.Range("B22").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image0'
$ws.Range("B24").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:22:Image1'
$ws.Range("B25").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image1'
$ws.Range("B29").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image0'
$ws.Range("B32").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image0'
$ws.Range("B33").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image1'
$ws.Range("B34").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image1'
$ws.Range("B36").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image0'
$ws.Range("B39").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:37:Image1'
$ws.Range("B42").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0'
$ws.Range("B43").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:41:Image1'
$ws.Range("B44").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:42:Image1'
$ws.Range("B45").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image0'
$ws.Range("B47").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:45:Image0'
$ws.Range("B52").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:50:Image0'
$ws.Range("B57").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:55:Image1'
$ws.Range("B59").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image0'
$ws.Range("B60").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:58:Image1'
$ws.Range("B61").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:59:Image1'
$ws.Range("B62").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:60:Image0'
$ws.Range("B65").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:63:Image0'
$ws.Range("B67").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:65:Image0'
$ws.Range("B69").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:67:Image0'
$ws.Range("B70").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image0'
$ws.Range("B73").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0'
$ws.Range("B77").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0'
$ws.Range("B78").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:76:Image1'
$ws.Range("B80").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:78:Image1'
$ws.Range("B81").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image1'

$ws = $wb.Worksheets.Item("SchemaDefinedTerm")
$ws.Range("B2").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:80:Image1'
$ws.Range("B5").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image0'
$ws.Range("B6").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image1'
$ws.Range("B7").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:85:Image1'
$ws.Range("B9").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image1'
$ws.Range("B10").Value = 'urn:paradicms_etl:pipeline:synthetic_data:concept:88:Image0'

# CreativeCommonsLicense: swap row 3 and row 4 contents (columns A-L)
$ws = $wb.Worksheets.Item("CreativeCommonsLicense")
$row3 = @()
$row4 = @()
for ($c = 1; $c -le 12; $c++) {
    $row3 += ,$ws.Cells.Item(3, $c).Value2
    $row4 += ,$ws.Cells.Item(4, $c).Value2
}
# Force text number format so numeric-looking strings (e.g. "1.0") round-trip
# as text instead of being coerced to a Double.
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(3, $c).NumberFormat = "@"
    $ws.Cells.Item(4, $c).NumberFormat = "@"
}
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(3, $c).Value2 = $row4[$c - 1]
    $ws.Cells.Item(4, $c).Value2 = $row3[$c - 1]
}
